$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 148.13333
$ws.Range("I9").Value = 148.13333
$ws.Range("K9").Value = 148.13333
$ws.Range("M9").Value = 20.86667
$ws.Range("H11").Value = 188.2
$ws.Range("I11").Value = 188.2
$ws.Range("K11").Value = 188.2
$ws.Range("M11").Value = -48.19999999999999
$ws.Range("H17").Value = 552.36365
$ws.Range("J17").Value = 341.77777
$ws.Range("L17").Value = 1025.33331
$ws.Range("N17").Value = -1361.33331
$ws.Range("H62").Value = 1503.3636
$ws.Range("I62").Value = 1503.3636
$ws.Range("K62").Value = 1503.3636
$ws.Range("M62").Value = -879.3635999999999
$ws.Range("H65").Value = 1503.3636
$ws.Range("I65").Value = 1503.3636
$ws.Range("K65").Value = 7516.817999999999
$ws.Range("M65").Value = -4396.817999999999
$ws.Range("H98").Value = 9357.556
$ws.Range("I98").Value = 8738.888999999999
$ws.Range("J98").Value = 9666.888999999999
$ws.Range("K98").Value = 8738.888999999999
$ws.Range("L98").Value = 9666.888999999999
$ws.Range("M98").Value = -7240.888999999999
$ws.Range("N98").Value = -12662.889
$ws.Range("H103").Value = 6434.933
$ws.Range("J103").Value = 18200
$ws.Range("L103").Value = 54600
$ws.Range("N103").Value = -55772
$ws.Range("H122").Value = 9357.556
$ws.Range("I122").Value = 8738.888999999999
$ws.Range("J122").Value = 9666.888999999999
$ws.Range("K122").Value = 26216.667
$ws.Range("L122").Value = 29000.667
$ws.Range("M122").Value = -23766.667
$ws.Range("N122").Value = -33900.667
$ws.Range("H137").Value = 1703045.1
$ws.Range("I137").Value = 2166111.8
$ws.Range("K137").Value = 6498335.399999999
$ws.Range("M137").Value = -6495785.399999999
$ws.Range("H138").Value = 2846.7827
$ws.Range("I138").Value = 1838.75
$ws.Range("K138").Value = 5516.25
$ws.Range("M138").Value = -376.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 682.2222
$ws.Range("I2").Value = 682.2222
$ws.Range("K2").Value = 682.2222
$ws.Range("M2").Value = -569.2222
$ws.Range("H30").Value = 7274.75
$ws.Range("I30").Value = 5600
$ws.Range("J30").Value = 8949.5
$ws.Range("K30").Value = 5600
$ws.Range("L30").Value = 8949.5
$ws.Range("M30").Value = -5450
$ws.Range("N30").Value = -9249.5
$ws.Range("H45").Value = 3254.4
$ws.Range("I45").Value = 3568
$ws.Range("K45").Value = 3568
$ws.Range("M45").Value = -3191
$ws.Range("H74").Value = 3018
$ws.Range("I74").Value = 2566.6667
$ws.Range("K74").Value = 2566.6667
$ws.Range("M74").Value = -1692.6667
$ws.Range("H77").Value = 3018
$ws.Range("I77").Value = 2566.6667
$ws.Range("K77").Value = 12833.3335
$ws.Range("M77").Value = -8465.333500000001
$ws.Range("H110").Value = 551.4
$ws.Range("I110").Value = 514
$ws.Range("J110").Value = 888
$ws.Range("K110").Value = 514
$ws.Range("L110").Value = 888
$ws.Range("M110").Value = 1531
$ws.Range("N110").Value = -4978
$ws.Range("H116").Value = 682.2222
$ws.Range("I116").Value = 682.2222
$ws.Range("K116").Value = 682.2222
$ws.Range("M116").Value = 1611.7778
$ws.Range("H132").Value = 2394.375
$ws.Range("I132").Value = 1081.7693
$ws.Range("K132").Value = 3245.3079
$ws.Range("M132").Value = -715.3078999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 682.2222
$ws.Range("I3").Value = 682.2222
$ws.Range("K3").Value = 682.2222
$ws.Range("M3").Value = -568.2222
$ws.Range("H20").Value = 24430
$ws.Range("I20").Value = 4755
$ws.Range("J20").Value = 32300
$ws.Range("K20").Value = 4755
$ws.Range("L20").Value = 32300
$ws.Range("M20").Value = -4508
$ws.Range("N20").Value = -32794
$ws.Range("H105").Value = 1787.5714
$ws.Range("I105").Value = 1798.091
$ws.Range("J105").Value = 1614
$ws.Range("K105").Value = 1798.091
$ws.Range("L105").Value = 1614
$ws.Range("M105").Value = -51.09099999999989
$ws.Range("N105").Value = -5108
$ws.Range("H107").Value = 1531.6154
$ws.Range("I107").Value = 1121.1
$ws.Range("J107").Value = 2900
$ws.Range("K107").Value = 1121.1
$ws.Range("L107").Value = 2900
$ws.Range("M107").Value = 798.9000000000001
$ws.Range("N107").Value = -6740

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4576.041
$ws.Range("I31").Value = 1759.75
$ws.Range("K31").Value = 1759.75
$ws.Range("M31").Value = -1464.75
$ws.Range("H34").Value = 4576.041
$ws.Range("I34").Value = 1759.75
$ws.Range("K34").Value = 1759.75
$ws.Range("M34").Value = -1557.75
$ws.Range("H58").Value = 1880.7354
$ws.Range("I58").Value = 1411.8966
$ws.Range("K58").Value = 1411.8966
$ws.Range("M58").Value = -1208.8966
$ws.Range("H132").Value = 3140.3333
$ws.Range("I132").Value = 1364.9231
$ws.Range("J132").Value = 7756.4
$ws.Range("K132").Value = 4094.7693
$ws.Range("L132").Value = 23269.2
$ws.Range("M132").Value = -1564.7693
$ws.Range("N132").Value = -28329.2
$ws.Range("H136").Value = 1880.7354
$ws.Range("I136").Value = 1411.8966
$ws.Range("K136").Value = 4235.6898
$ws.Range("M136").Value = -1685.6898

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 6326.3335
$ws.Range("J62").Value = 8989.5
$ws.Range("L62").Value = 26968.5
$ws.Range("N62").Value = -28340.5
$ws.Range("H65").Value = 6326.3335
$ws.Range("J65").Value = 8989.5
$ws.Range("L65").Value = 80905.5
$ws.Range("N65").Value = -87769.5
$ws.Range("H122").Value = 2897.625
$ws.Range("J122").Value = 3248.9
$ws.Range("L122").Value = 29240.1
$ws.Range("N122").Value = -34140.10000000001
$ws.Range("H131").Value = 786.1
$ws.Range("J131").Value = 822.25806
$ws.Range("L131").Value = 2466.77418
$ws.Range("N131").Value = -12546.77418

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1861.409
$ws.Range("I113").Value = 1968.1428
$ws.Range("J113").Value = 1674.625
$ws.Range("K113").Value = 1968.1428
$ws.Range("L113").Value = 1674.625
$ws.Range("M113").Value = 201.8571999999999
$ws.Range("N113").Value = -6014.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 12458
$ws.Range("J29").Value = 20000
$ws.Range("L29").Value = 20000
$ws.Range("N29").Value = -20590
$ws.Range("H46").Value = 1812.9412
$ws.Range("I46").Value = 1247.2727
$ws.Range("J46").Value = 2850
$ws.Range("K46").Value = 1247.2727
$ws.Range("L46").Value = 2850
$ws.Range("M46").Value = -1059.2727
$ws.Range("N46").Value = -3226
$ws.Range("H132").Value = 8588.111000000001
$ws.Range("I132").Value = 3548.5
$ws.Range("J132").Value = 12619.8
$ws.Range("K132").Value = 10645.5
$ws.Range("L132").Value = 37859.39999999999
$ws.Range("M132").Value = -8115.5
$ws.Range("N132").Value = -42919.39999999999
$ws.Range("H136").Value = 5047.8423
$ws.Range("I136").Value = 1587
$ws.Range("K136").Value = 4761
$ws.Range("M136").Value = -2211

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 49400.9
$ws.Range("J46").Value = 49400.9
$ws.Range("L46").Value = 49400.9
$ws.Range("N46").Value = -49862.9
$ws.Range("H134").Value = 49400.9
$ws.Range("J134").Value = 49400.9
$ws.Range("L134").Value = 148202.7
$ws.Range("N134").Value = -153272.7
$ws.Range("H136").Value = 7332.533
$ws.Range("I136").Value = 7818.294
$ws.Range("J136").Value = 6697.3076
$ws.Range("K136").Value = 23454.882
$ws.Range("L136").Value = 20091.9228
$ws.Range("M136").Value = -20904.882
$ws.Range("N136").Value = -25191.9228
